$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates: column, row, new text value.
# NumberFormat is forced to text ("@") before writing so Excel does not
# auto-convert these numeric-looking / percent-looking strings into
# actual numbers (the source file stores them as plain text).
$updates = @(
    @{ Cell = 'D2'; Value = '306.47' }
    @{ Cell = 'E2'; Value = '-3.37%' }
    @{ Cell = 'D3'; Value = '37.47' }
    @{ Cell = 'E3'; Value = '-5.99%' }
    @{ Cell = 'D4'; Value = '5.083' }
    @{ Cell = 'E4'; Value = '-1.37%' }
    @{ Cell = 'D5'; Value = '0.07729' }
    @{ Cell = 'E5'; Value = '-6.23%' }
    @{ Cell = 'D6'; Value = '4.335' }
    @{ Cell = 'E6'; Value = '0.14%' }
    @{ Cell = 'D7'; Value = '1.884' }
    @{ Cell = 'E7'; Value = '-8.16%' }
    @{ Cell = 'D8'; Value = '8.178' }
    @{ Cell = 'E8'; Value = '-2.25%' }
    @{ Cell = 'E9'; Value = '-8.38%' }
    @{ Cell = 'D10'; Value = '0.9236' }
    @{ Cell = 'E10'; Value = '-1.73%' }
    @{ Cell = 'D11'; Value = '0.1231' }
    @{ Cell = 'E11'; Value = '-8.59%' }
    @{ Cell = 'D12'; Value = '0.1870' }
    @{ Cell = 'E12'; Value = '-6.33%' }
    @{ Cell = 'D13'; Value = '0.08764' }
    @{ Cell = 'E13'; Value = '-3.12%' }
    @{ Cell = 'D14'; Value = '0.03405' }
    @{ Cell = 'E14'; Value = '-3.37%' }
    @{ Cell = 'D15'; Value = '0.09706' }
    @{ Cell = 'E15'; Value = '-0.94%' }
    @{ Cell = 'D16'; Value = '0.001370' }
    @{ Cell = 'E16'; Value = '-2.90%' }
    @{ Cell = 'D17'; Value = '0.006077' }
    @{ Cell = 'E17'; Value = '0.66%' }
    @{ Cell = 'D18'; Value = '3.586' }
    @{ Cell = 'E18'; Value = '-2.60%' }
    @{ Cell = 'E19'; Value = '-2.46%' }
    @{ Cell = 'D20'; Value = '0.1269' }
    @{ Cell = 'E20'; Value = '-4.16%' }
    @{ Cell = 'D21'; Value = '5.016' }
    @{ Cell = 'E22'; Value = '1.50%' }
    @{ Cell = 'D23'; Value = '0.02104' }
    @{ Cell = 'E23'; Value = '5,157.19%' }
    @{ Cell = 'D24'; Value = '0.04329' }
    @{ Cell = 'E24'; Value = '-0.61%' }
    @{ Cell = 'E25'; Value = '-2.23%' }
    @{ Cell = 'D26'; Value = '0.004229' }
    @{ Cell = 'E26'; Value = '-11.91%' }
    @{ Cell = 'E27'; Value = '3.68%' }
    @{ Cell = 'D39'; Value = '0.02180' }
    @{ Cell = 'E39'; Value = '-5.26%' }
    @{ Cell = 'E40'; Value = '-5.64%' }
    @{ Cell = 'D41'; Value = '0.007484' }
    @{ Cell = 'E41'; Value = '-3.10%' }
    @{ Cell = 'D42'; Value = '0.009883' }
    @{ Cell = 'E42'; Value = '-5.02%' }
    @{ Cell = 'D43'; Value = '0.1339' }
    @{ Cell = 'E43'; Value = '-4.92%' }
    @{ Cell = 'E44'; Value = '-1.02%' }
    @{ Cell = 'D45'; Value = '0.009836' }
    @{ Cell = 'E45'; Value = '5.53%' }
    @{ Cell = 'D46'; Value = '0.00006536' }
    @{ Cell = 'E46'; Value = '-9.70%' }
    @{ Cell = 'E47'; Value = '-0.07%' }
    @{ Cell = 'E48'; Value = '3.79%' }
    @{ Cell = 'E50'; Value = '-0.07%' }
    @{ Cell = 'E51'; Value = '-0.07%' }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    $range.NumberFormat = "@"
    $range.Value = $u.Value
}
